# Insert two new weekly price rows before the existing row 110
# (Agricola del Norte S.A. de Arica - Caigua, week of 44841) and shift
# the remaining rows (former 110-139) down to 112-141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 110 - pushes old rows 110..139 down to 112..141
$ws.Rows.Item(110).Resize(2).Insert()

# Row 110: Primera
$ws.Cells.Item(110, 1).Value = 1
$ws.Cells.Item(110, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(110, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(110, 4).Value = 44841
$ws.Cells.Item(110, 5).Value = 15
$ws.Cells.Item(110, 6).Value = 100112036
$ws.Cells.Item(110, 7).Value = "Caigua"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 130
$ws.Cells.Item(110, 11).Value = 7000
$ws.Cells.Item(110, 12).Value = 8000
$ws.Cells.Item(110, 13).Value = 7500
$ws.Cells.Item(110, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(110, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(110, 16).Value = 375
$ws.Cells.Item(110, 17).Value = 20
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Row 111: Segunda
$ws.Cells.Item(111, 1).Value = 1
$ws.Cells.Item(111, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(111, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(111, 4).Value = 44841
$ws.Cells.Item(111, 5).Value = 15
$ws.Cells.Item(111, 6).Value = 100112036
$ws.Cells.Item(111, 7).Value = "Caigua"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Segunda"
$ws.Cells.Item(111, 10).Value = 140
$ws.Cells.Item(111, 11).Value = 6000
$ws.Cells.Item(111, 12).Value = 7000
$ws.Cells.Item(111, 13).Value = 6500
$ws.Cells.Item(111, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(111, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(111, 16).Value = 325
$ws.Cells.Item(111, 17).Value = 20
$ws.Cells.Item(111, 18).Value = "Hortaliza"

# Match the date formatting used by the rest of column D (style applied via
# the row insert should already copy it from row 109, but set explicitly too)
$ws.Range("D110:D111").NumberFormat = $ws.Range("D109").NumberFormat
